$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update E8: BSides Mumbai CTF 2024 description (sponsor shoutout + discord bot artifact) ---
$ws.Cells.Item(8, 5).Value = 'Join us for the inaugural BSides Mumbai Capture The Flag (CTF) competition! This online competition, a Jeopardy Style CTF, is designed for beginners, intermediate players, and all ethical hackers. All are invited to participate! So get ready to unleash your true potential and discover the hacker in you. The competition will feature a series of challenges spanning cryptography, reverse engineering, web security, and more, catering to participants of all skill levels.

We have lot of Categories :
1. Web
2. Forensic
3. Osint
4. Crypto
5. Network
6. Cloud
7. Miscellaneous

Contact:
For questions or more information, contact
info@bsidesmumbai.in
Special Thanks to our Sponsors:
HackTheBox
Altered Security
Follow us on social media:
- LinkedIn: https://in.linkedin.com/company/bsidesmumbai
- Twitter: https://twitter.com/BSidesMumbai
- Instagram: https://www.instagram.com/bsidesmumbai/
- Discord: https://discord.gg/2KRGQWBGR3
Spea'

# --- Append new CTF event rows 24-34 ---
# Row 24
$ws.Cells.Item(24, 1).Value = 2304
$ws.Cells.Item(24, 2).Value = 'idekCTF 2024'
$ws.Cells.Item(24, 3).Value = '2024-06-15T00:00:00+00:00'
$ws.Cells.Item(24, 4).Value = '2024-06-17T00:00:00+00:00'
$ws.Cells.Item(24, 5).Value = 'idekCTF is an information security CTF competition organized by the idek team and is aimed at the entire spectrum from high school and university students to experienced players. idekCTF will cover the standard Jeopardy-style CTF topics (binary exploitation, reverse engineering, cryptography, web exploitation, and forensics) as well as other, less standard categories. '
$ws.Cells.Item(24, 6).Value = 'https://ctf.idek.team/'
$ws.Cells.Item(24, 7).Value = 'idek'
$ws.Cells.Item(24, 8).Value = 48
$ws.Cells.Item(24, 9).Value = $true
$ws.Cells.Item(24, 10).Value = $false
$ws.Cells.Item(24, 11).Value = $false

# Row 25
$ws.Cells.Item(25, 1).Value = 2342
$ws.Cells.Item(25, 2).Value = 'justCTF 2024'
$ws.Cells.Item(25, 3).Value = '2024-06-15T08:00:00+00:00'
$ws.Cells.Item(25, 4).Value = '2024-06-16T08:00:00+00:00'
$ws.Cells.Item(25, 5).Value = 'Sponsors: TBA'
$ws.Cells.Item(25, 6).Value = 'http://2024.justctf.team/'
$ws.Cells.Item(25, 7).Value = 'justCatTheFish'
$ws.Cells.Item(25, 8).Value = 24
$ws.Cells.Item(25, 9).Value = $true
$ws.Cells.Item(25, 10).Value = $false
$ws.Cells.Item(25, 11).Value = $false

# Row 26
$ws.Cells.Item(26, 1).Value = 2296
$ws.Cells.Item(26, 2).Value = 'Google Capture The Flag 2024'
$ws.Cells.Item(26, 3).Value = '2024-06-21T18:00:00+00:00'
$ws.Cells.Item(26, 4).Value = '2024-06-23T18:00:00+00:00'
$ws.Cells.Item(26, 5).Value = 'categories: web, pwn, crypto, sandbox, reversing, misc'
$ws.Cells.Item(26, 6).Value = 'https://g.co/ctf'
$ws.Cells.Item(26, 7).Value = 'Google CTF'
$ws.Cells.Item(26, 8).Value = 48
$ws.Cells.Item(26, 9).Value = $true
$ws.Cells.Item(26, 10).Value = $false
$ws.Cells.Item(26, 11).Value = $false

# Row 27
$ws.Cells.Item(27, 1).Value = 2275
$ws.Cells.Item(27, 2).Value = 'UIUCTF 2024'
$ws.Cells.Item(27, 3).Value = '2024-06-29T00:00:00+00:00'
$ws.Cells.Item(27, 4).Value = '2024-07-01T00:00:00+00:00'
$ws.Cells.Item(27, 5).Value = 'UIUCTF will be returning in 2024!'
$ws.Cells.Item(27, 6).Value = 'https://uiuc.tf/'
$ws.Cells.Item(27, 7).Value = 'SIGPwny'
$ws.Cells.Item(27, 8).Value = 48
$ws.Cells.Item(27, 9).Value = $true
$ws.Cells.Item(27, 10).Value = $false
$ws.Cells.Item(27, 11).Value = $false

# Row 28
$ws.Cells.Item(28, 1).Value = 2259
$ws.Cells.Item(28, 2).Value = 'Junior.Crypt.2024 CTF'
$ws.Cells.Item(28, 3).Value = '2024-07-03T15:00:00+00:00'
$ws.Cells.Item(28, 4).Value = '2024-07-05T15:00:00+00:00'
$ws.Cells.Item(28, 5).Value = 'Junior.Crypt.2024 CTF is an open competition in information security for beginners, students and everyone. This is a 48 hour online event. We hope all participants can use their skills and have a great time.
The organizers of the CTF were students of the Department of System Programming and Computer Security of Grodno State University, Grodno, Belarus.
Team participation. Team of 1-4 people.
The tasks relate to different areas - beginner, cryptography, forensics, miscellaneous, OSINT, ppc, pwn, rev, web. They will be of interest to beginners, students and everyone else.
The official languages of the tournament are English and Russian. But you can take part even if you don''t speak any of them. We will do our best to ensure that each condition of the problem can be understood with the correct use of machine translation. In addition, during the tournament tasks, we will try to introduce our participants to an unusually beautiful country, our Republic of Belarus.'
$ws.Cells.Item(28, 6).Value = 'http://ctf-spcs.mf.grsu.by/'
$ws.Cells.Item(28, 7).Value = 'Beavers0'
$ws.Cells.Item(28, 8).Value = 48
$ws.Cells.Item(28, 9).Value = $true
$ws.Cells.Item(28, 10).Value = $false
$ws.Cells.Item(28, 11).Value = $false

# Row 29
$ws.Cells.Item(29, 1).Value = 2284
$ws.Cells.Item(29, 2).Value = 'DownUnderCTF 2024'
$ws.Cells.Item(29, 3).Value = '2024-07-05T09:30:00+00:00'
$ws.Cells.Item(29, 4).Value = '2024-07-07T09:30:00+00:00'
$ws.Cells.Item(29, 5).Value = 'DownUnderCTF is the largest online Australian-run Capture The Flag (CTF) competition, now welcoming Aotearoa (New Zealand) to the competition for the first time in 2024. With over 4200+ registered users and more than 2000+ registered teams as of 2023, its primary goal is to up-skill the next generation of potential Cyber Security Professionals and to expand the CTF community in Australia and Aotearoa (New Zealand). While our CTF is an online event open to participants worldwide, starting from 2024, prize eligibility extends to include both Australian and Aotearoa (New Zealand) Secondary or Tertiary school students. This change aims to foster a closer collaboration and competition spirit between the two nations while maintaining our commitment to enhancing cybersecurity skills among the youth.'
$ws.Cells.Item(29, 6).Value = 'https://play.duc.tf/'
$ws.Cells.Item(29, 7).Value = 'DownUnderCTF'
$ws.Cells.Item(29, 8).Value = 48
$ws.Cells.Item(29, 9).Value = $true
$ws.Cells.Item(29, 10).Value = $false
$ws.Cells.Item(29, 11).Value = $false

# Row 30
$ws.Cells.Item(30, 1).Value = 2301
$ws.Cells.Item(30, 2).Value = 'Interlogica CTF2024 - Wastelands'
$ws.Cells.Item(30, 3).Value = '2024-07-05T12:37:00+00:00'
$ws.Cells.Item(30, 4).Value = '2024-07-07T22:59:59+00:00'
$ws.Cells.Item(30, 5).Value = 'Welcome to wastelands, where the faint echoes of civilization cling to the sands like whispers of a bygone era. Abandoned technology and rusting drones dot the barren landscape, serving as eerie reminders of humanity''s downfall. Amidst the desolation, the remnants of once-thriving sustainability systems cling desperately to survival, struggling to eke out existence amidst the chaos.
Prepare to confront the harsh realities of this unforgiving terrain as you embark on thrilling challenges that will push you to your limits: Reviving Relics, Hydrophonic Systems, Unlocking the Vault, Digital Infiltration, Satellite Takeover, Expeditions, Vehicle Restoration are some of the challenges that await you!
Embark on this epic odyssey and prove your mettle in a world where survival demands more than just strength of arms—it demands cunning, resilience, and the indomitable will to endure. Join us in the ultimate test of endurance and strategy amidst the ruins of civilization. Will you rise to the challenge, or be swallowed by the sands of time?
Event Start:  2024 July 5 12:37
Event End: 2024 July 7 22:59
Maximum Team size: 5 members
Rulez:
1. During brute-force attacks on services and/or web pages, the use of the specified wordlists is mandatory.
2. Actions that overload the resources of competition services, such as saturating disk space or CPU, are prohibited. These activities will be constantly monitored.
3. Exchange of flags and/or solutions between different teams is strictly prohibited. Publishing of walkthroughs or flags before the end of the event is prohibited. Once the event has ended, publishing is encouraged.
4. Performing attacks on infrastructures other than those specified by the challenges is prohibited.
5. The use of multiple environments and different types of automatic encoding presupposes that the participant is able to handle them independently. In case of uncertainty, it is recommended to use a Kali/Parrot virtual machine.
6. Each violation will result in a loss of points for the entire team.
7. The team that first solves a specific challenge will receive additional points (first blood).
8. The top 3 ranked teams will be awarded.
9. Registrations on our platform will open one month before the event.
10. Teams can consist of up to a maximum of 5 participants.
11. If an individual does not have a team, a dedicated section for team formation will be available on Discord.
12. Each team will be provided with a dedicated WireGuard VPN to tackle the challenges. '
$ws.Cells.Item(30, 6).Value = 'https://ctf.interlogica.it/'
$ws.Cells.Item(30, 7).Value = 'Interlogica'
$ws.Cells.Item(30, 8).Value = 58
$ws.Cells.Item(30, 9).Value = $true
$ws.Cells.Item(30, 10).Value = $false
$ws.Cells.Item(30, 11).Value = $false

# Row 31
$ws.Cells.Item(31, 1).Value = 2345
$ws.Cells.Item(31, 2).Value = 'HITCON CTF 2024 Quals'
$ws.Cells.Item(31, 3).Value = '2024-07-12T14:00:00+00:00'
$ws.Cells.Item(31, 4).Value = '2024-07-14T14:00:00+00:00'
$ws.Cells.Item(31, 5).Value = 'TBA'
$ws.Cells.Item(31, 6).Value = 'http://ctf.hitcon.org/'
$ws.Cells.Item(31, 7).Value = 'HITCON'
$ws.Cells.Item(31, 8).Value = 48
$ws.Cells.Item(31, 9).Value = $true
$ws.Cells.Item(31, 10).Value = $false
$ws.Cells.Item(31, 11).Value = $false

# Row 32
$ws.Cells.Item(32, 1).Value = 2293
$ws.Cells.Item(32, 2).Value = 'MOCA CTF - Qualification'
$ws.Cells.Item(32, 3).Value = '2024-07-20T09:00:00+00:00'
$ws.Cells.Item(32, 4).Value = '2024-07-21T09:00:00+00:00'
$ws.Cells.Item(32, 5).Value = '[DELAYED TO 20-21 July]
Official CTF competition of the Metro Olografix Camp, organized by MOCA, Fibonhack and PWNX.
Best teams will be invited to compete at the final events with travel expenses reinbursement!'
$ws.Cells.Item(32, 6).Value = 'https://moca.camp/ctf/'
$ws.Cells.Item(32, 7).Value = 'Metro Olografix'
$ws.Cells.Item(32, 8).Value = 24
$ws.Cells.Item(32, 9).Value = $true
$ws.Cells.Item(32, 10).Value = $false
$ws.Cells.Item(32, 11).Value = $false

# Row 33
$ws.Cells.Item(33, 1).Value = 2353
$ws.Cells.Item(33, 2).Value = 'DeadSec CTF 2024'
$ws.Cells.Item(33, 3).Value = '2024-07-26T20:00:00+00:00'
$ws.Cells.Item(33, 4).Value = '2024-07-28T08:00:00+00:00'
$ws.Cells.Item(33, 5).Value = 'DeadSec CTF 2024 is an online jeopardy-style CTF organized by DeadSec Team.
There will be challenges with a wide range of difficulty mainly from cryptography, reverse, pwn, web, misc...'
# F33 is empty in source (no value assigned)
$ws.Cells.Item(33, 7).Value = 'DeadSec'
$ws.Cells.Item(33, 8).Value = 36
$ws.Cells.Item(33, 9).Value = $true
$ws.Cells.Item(33, 10).Value = $false
$ws.Cells.Item(33, 11).Value = $false

# Row 34
$ws.Cells.Item(34, 1).Value = 2282
$ws.Cells.Item(34, 2).Value = 'corCTF 2024'
$ws.Cells.Item(34, 3).Value = '2024-07-27T00:00:00+00:00'
$ws.Cells.Item(34, 4).Value = '2024-07-29T00:00:00+00:00'
$ws.Cells.Item(34, 5).Value = 'Infra sponsored by <a href="https://goo.gle/ctfsponsorship">goo.gle/ctfsponsorship</a>'
$ws.Cells.Item(34, 6).Value = 'https://ctf.cor.team/'
$ws.Cells.Item(34, 7).Value = 'Crusaders of Rust'
$ws.Cells.Item(34, 8).Value = 48
$ws.Cells.Item(34, 9).Value = $true
$ws.Cells.Item(34, 10).Value = $false
$ws.Cells.Item(34, 11).Value = $false

Write-Output "done"